$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "0.535", "1.99", "0.140") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# and the original literal formatting (e.g. trailing zeros) would be lost.
$textCells = @(
    "D5",
    "D6",
    "D9",
    "D12",
    "D19",
    "D20",
    "D22",
    "D24",
    "D28",
    "D34",
    "D37",
    "D38",
    "D39",
    "D42",
    "D43",
    "D44",
    "D45",
    "D49",
    "D50"
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated cryptocurrency price / link / volume values
$ws.Range("D2").Value = '69.152.68'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '3.813.68'
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '601.84'
$ws.Range("E5").Value = '  -0.73%  '
$ws.Range("D6").Value = '164.12'
$ws.Range("E6").Value = '  -3.52%  '
$ws.Range("D7").Value = '3.811.14'
$ws.Range("E7").Value = '  +1.14%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.535'
$ws.Range("E10").Value = '  +0.66%  '
$ws.Range("E11").Value = '  -1.80%  '
$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("E14").Value = '  -2.24%  '
$ws.Range("D15").Value = '4.451.93'
$ws.Range("E15").Value = '  +1.20%  '
$ws.Range("D16").Value = '3.812.89'
$ws.Range("E16").Value = '  +1.17%  '
$ws.Range("D17").Value = '69.277.30'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E18").Value = '  +1.50%  '
$ws.Range("D19").Value = '11.54'
$ws.Range("E19").Value = '  +5.65%  '
$ws.Range("D20").Value = '17.39'
$ws.Range("E20").Value = '  +1.17%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").Value = '489.56'
$ws.Range("E22").Value = '  -1.63%  '
$ws.Range("E23").Value = '  -1.20%  '
$ws.Range("D24").Value = '0.0000157'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("E25").Value = '  -0.89%  '
$ws.Range("E26").Value = '  -3.92%  '
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").Value = '10.05'
$ws.Range("E28").Value = '  -2.99%  '
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("E32").Value = '  -4.65%  '
$ws.Range("D33").Value = '3.961.38'
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("D34").Value = '31.97'
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("D35").Value = '3.760.94'
$ws.Range("E35").Value = '  +1.53%  '
$ws.Range("E36").Value = '  -1.98%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.140'
$ws.Range("E37").Value = '  +4.76%  '
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = '5.93'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").Value = '  -0.01%  '
$ws.Range("E41").Value = '  -1.87%  '
$ws.Range("D42").Value = '3.02'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("D43").Value = '48.62'
$ws.Range("E43").Value = '  -0.25%  '
$ws.Range("D44").Value = '1.99'
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").Value = '425.82'
$ws.Range("E45").Value = '  -3.82%  '
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = '2.833.18'
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").Value = '141.35'
$ws.Range("E49").Value = '  +0.36%  '
$ws.Range("D50").Value = '39.41'
$ws.Range("E50").Value = '  -3.48%  '
$ws.Range("E51").Value = '  -1.74%  '
